$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.836.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.735.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.37%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5181"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2739"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06147"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.740.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07165"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("E12").Value = "  -2.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6398"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.603"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.13"
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.874.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("E19").Value = "  +1.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006762"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.623"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.247"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.518"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.767"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.945"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08249"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.654"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04628"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.643"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9851"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.70%  "

$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.73%  "

$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.916"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3836"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7465"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.50%  "

$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1122"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.230"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("E47").Value = "  -1.83%  "

$ws.Range("E48").Value = "  +1.98%  "

$ws.Range("E49").Value = "  +1.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.556"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3405"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.74%  "
